$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2387.5
$ws.Range("I6").Value = 2350
$ws.Range("J6").Value = 2500
$ws.Range("K6").Value = 7050
$ws.Range("L6").Value = 7500
$ws.Range("M6").Value = -6938
$ws.Range("N6").Value = -7724
$ws.Range("H8").Value = 111111460
$ws.Range("H38").Value = 689.0952
$ws.Range("I38").Value = 98.066666
$ws.Range("J38").Value = 2166.6667
$ws.Range("K38").Value = 294.199998
$ws.Range("L38").Value = 6500.000100000001
$ws.Range("M38").Value = 77.80000200000001
$ws.Range("N38").Value = -7244.000100000001
$ws.Range("H39").Value = 142.15384
$ws.Range("I39").Value = 32.8
$ws.Range("J39").Value = 506.66666
$ws.Range("K39").Value = 98.39999999999999
$ws.Range("L39").Value = 1519.99998
$ws.Range("M39").Value = 197.6
$ws.Range("N39").Value = -2111.99998
$ws.Range("H42").Value = 175.2
$ws.Range("I42").Value = 91.5
$ws.Range("J42").Value = 510
$ws.Range("K42").Value = 274.5
$ws.Range("L42").Value = 1530
$ws.Range("M42").Value = -44.5
$ws.Range("N42").Value = -1990
$ws.Range("H58").Value = 1642.75
$ws.Range("I58").Value = 357.5
$ws.Range("J58").Value = 5498.5
$ws.Range("K58").Value = 1072.5
$ws.Range("L58").Value = 16495.5
$ws.Range("M58").Value = -922.5
$ws.Range("N58").Value = -16795.5
$ws.Range("H87").Value = 28516.666
$ws.Range("J87").Value = 28516.666
$ws.Range("L87").Value = 28516.666
$ws.Range("N87").Value = -31012.666
$ws.Range("H90").Value = 28516.666
$ws.Range("J90").Value = 28516.666
$ws.Range("L90").Value = 85549.99800000001
$ws.Range("N90").Value = -98029.99800000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 276.33334
$ws.Range("I4").Value = 276.33334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 276.33334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -160.33334
$ws.Range("H5").Value = 75.71429000000001
$ws.Range("I5").Value = 71.666664
$ws.Range("K5").Value = 71.666664
$ws.Range("M5").Value = 40.333336
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = ""
$ws.Range("N12").Value = 0
$ws.Range("H34").Value = 12375
$ws.Range("J34").Value = 12375
$ws.Range("L34").Value = 12375
$ws.Range("N34").Value = -12917
$ws.Range("H55").Value = 9833.333000000001
$ws.Range("J55").Value = 9833.333000000001
$ws.Range("L55").Value = 9833.333000000001
$ws.Range("N55").Value = -10463.333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 75.71429000000001
$ws.Range("I4").Value = 71.666664
$ws.Range("K4").Value = 71.666664
$ws.Range("M4").Value = 43.333336
$ws.Range("H7").Value = 2108.3
$ws.Range("I7").Value = 730.5
$ws.Range("J7").Value = 4175
$ws.Range("K7").Value = 730.5
$ws.Range("L7").Value = 4175
$ws.Range("M7").Value = -617.5
$ws.Range("N7").Value = -4401
$ws.Range("H76").Value = 26530.8
$ws.Range("J76").Value = 26530.8
$ws.Range("L76").Value = 26530.8
$ws.Range("N76").Value = -27160.8
$ws.Range("H79").Value = 26530.8
$ws.Range("J79").Value = 26530.8
$ws.Range("L79").Value = 26530.8
$ws.Range("N79").Value = -28714.8
$ws.Range("H80").Value = 755.625
$ws.Range("I80").Value = 299
$ws.Range("J80").Value = 907.8333
$ws.Range("K80").Value = 299
$ws.Range("L80").Value = 907.8333
$ws.Range("M80").Value = 699
$ws.Range("N80").Value = -2903.8333
$ws.Range("H83").Value = 755.625
$ws.Range("I83").Value = 299
$ws.Range("J83").Value = 907.8333
$ws.Range("K83").Value = 1495
$ws.Range("L83").Value = 4539.1665
$ws.Range("M83").Value = 3497
$ws.Range("N83").Value = -14523.1665
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 501
$ws.Range("I12").Value = 501
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 501
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = -331
$ws.Range("H35").Value = 200000980
$ws.Range("I35").Value = 200000980
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 200000980
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = -200000686
$ws.Range("H50").Value = 11500
$ws.Range("I50").Value = 3000
$ws.Range("J50").Value = 20000
$ws.Range("K50").Value = 3000
$ws.Range("L50").Value = 20000
$ws.Range("M50").Value = -2375
$ws.Range("N50").Value = -21250
$ws.Range("H99").Value = 38465.75
$ws.Range("I99").Value = 93410.37
$ws.Range("J99").Value = 2913.353
$ws.Range("K99").Value = 93410.37
$ws.Range("L99").Value = 2913.353
$ws.Range("M99").Value = -91912.37
$ws.Range("N99").Value = -5909.353
$ws.Range("H126").Value = 38465.75
$ws.Range("I126").Value = 93410.37
$ws.Range("J126").Value = 2913.353
$ws.Range("K126").Value = 280231.11
$ws.Range("L126").Value = 8740.059000000001
$ws.Range("M126").Value = -277761.11
$ws.Range("N126").Value = -13680.059
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 824.58826
$ws.Range("I34").Value = 951.8
$ws.Range("J34").Value = 771.5833
$ws.Range("K34").Value = 2855.4
$ws.Range("L34").Value = 2314.7499
$ws.Range("M34").Value = -2771.4
$ws.Range("N34").Value = -2482.7499
$ws.Range("H39").Value = 4166.6665
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4166.6665
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = ""
$ws.Range("M39").Value = 12499.9995
$ws.Range("N39").Value = -13087.9995
$ws.Range("H55").Value = 2922.125
$ws.Range("J55").Value = 3633.3333
$ws.Range("L55").Value = 10899.9999
$ws.Range("N55").Value = -11253.9999
$ws.Range("H59").Value = 1947.5
$ws.Range("I59").Value = 990
$ws.Range("J59").Value = 2266.6667
$ws.Range("K59").Value = 2970
$ws.Range("L59").Value = 6800.000100000001
$ws.Range("M59").Value = -2430
$ws.Range("N59").Value = -7880.000100000001
$ws.Range("H70").Value = 3531.8572
$ws.Range("I70").Value = 907.6667
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 2723.0001
$ws.Range("L70").Value = 16500
$ws.Range("M70").Value = -2408.0001
$ws.Range("N70").Value = -17130
$ws.Range("H73").Value = 3531.8572
$ws.Range("I73").Value = 907.6667
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 2723.0001
$ws.Range("L73").Value = 16500
$ws.Range("M73").Value = -1631.0001
$ws.Range("N73").Value = -18684
$ws.Range("H131").Value = 792.66174
$ws.Range("I131").Value = 325.23077
$ws.Range("J131").Value = 903.14545
$ws.Range("K131").Value = 975.69231
$ws.Range("L131").Value = 2709.43635
$ws.Range("M131").Value = 4064.30769
$ws.Range("N131").Value = -12789.43635
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 50000000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""
$ws.Range("H8").Value = 50000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = ""
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2343.611
$ws.Range("I7").Value = 2057.9167
$ws.Range("K7").Value = 2057.9167
$ws.Range("M7").Value = -1945.9167
$ws.Range("H126").Value = 2343.611
$ws.Range("I126").Value = 2057.9167
$ws.Range("K126").Value = 6173.750100000001
$ws.Range("M126").Value = -3703.750100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -858
$ws.Range("H23").Value = 265
$ws.Range("I23").Value = 265
$ws.Range("K23").Value = 265
$ws.Range("M23").Value = -36
$ws.Range("H68").Value = 35831
$ws.Range("J68").Value = 35831
$ws.Range("L68").Value = 35831
$ws.Range("N68").Value = -37453
$ws.Range("H71").Value = 35831
$ws.Range("J71").Value = 35831
$ws.Range("L71").Value = 107493
$ws.Range("N71").Value = -115605

Write-Output "edits applied"